# Applies the GitHub Actions crypto-price refresh (Mon Oct 30 00:45:14 UTC 2023).
# D-column price strings are number-like text (e.g. "34.516.38", "0.00...")
# so they are written with a leading quote-prefix to force Excel to keep them
# as literal text instead of silently re-parsing/truncating them as numbers
# (this mirrors typing '4.00 directly into a cell in the Excel UI).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.516.38"
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").Value = "'1.793.23"
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = "'226.82"
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").Value = "'0.555"
$ws.Range("E6").Value = '  +1.63%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = "'32.63"
$ws.Range("E8").Value = '  +2.82%  '

$ws.Range("E9").Value = '  +0.77%  '

$ws.Range("D10").Value = "'0.0691"
$ws.Range("E10").Value = '  +0.12%  '

$ws.Range("D11").Value = "'0.0948"
$ws.Range("E11").Value = '  +0.11%  '

$ws.Range("D12").Value = "'2.053.70"
$ws.Range("E12").Value = '  +0.82%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = "'11.09"
$ws.Range("E13").Value = '  +1.49%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = "'1.784.88"
$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").Value = "'0.636"
$ws.Range("E15").Value = '  +2.39%  '

$ws.Range("D16").Value = "'34.522.12"
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D18").Value = "'68.80"
$ws.Range("E18").Value = '  +1.32%  '

$ws.Range("D19").Value = "'0.0₃0800"
$ws.Range("E19").Value = '  +0.84%  '

$ws.Range("D20").Value = "'245.79"
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").Value = "'11.35"
$ws.Range("E21").Value = '  +2.63%  '

$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").Value = "'4.15"
$ws.Range("E23").Value = '  +1.06%  '

$ws.Range("D24").Value = "'173.72"
$ws.Range("E24").Value = '  +7.52%  '

$ws.Range("D25").Value = "'2.05"
$ws.Range("E25").Value = '  +0.43%  '

$ws.Range("E26").Value = '  +2.06%  '

$ws.Range("D27").Value = "'16.57"
$ws.Range("E27").Value = '  +1.78%  '

$ws.Range("D28").Value = "'0.115"
$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("D30").Value = "'4.00"
$ws.Range("E30").Value = '  +7.59%  '

$ws.Range("E31").Value = '  +0.71%  '

$ws.Range("D32").Value = "'0.0523"
$ws.Range("E32").Value = '  +0.59%  '

$ws.Range("D33").Value = "'3.78"
$ws.Range("E33").Value = '  +1.66%  '

$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = '  +1.44%  '

$ws.Range("D35").Value = "'1.421.75"
$ws.Range("E35").Value = '  -1.43%  '

$ws.Range("D36").Value = "'2.55"
$ws.Range("E36").Value = '  +5.68%  '

$ws.Range("D37").Value = "'0.674"
$ws.Range("E37").Value = '  +2.64%  '

$ws.Range("E38").Value = '  +2.56%  '

$ws.Range("D39").Value = "'0.0190"
$ws.Range("E39").Value = '  -0.26%  '

$ws.Range("D40").Value = "'84.35"
$ws.Range("E40").Value = '  +5.03%  '

$ws.Range("D41").Value = "'0.943"
$ws.Range("E41").Value = '  +1.84%  '

$ws.Range("D42").Value = "'2.38"
$ws.Range("E42").Value = '  +0.94%  '

$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = '  +2.84%  '

$ws.Range("D44").Value = "'13.93"
$ws.Range("E44").Value = '  +3.98%  '

$ws.Range("D45").Value = "'0.0527"
$ws.Range("E45").Value = '  +3.36%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = "'1.09"
$ws.Range("E46").Value = '  +1.87%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'6.10"
$ws.Range("E47").Value = '  +0.35%  '

$ws.Range("D48").Value = "'1.954.75"
$ws.Range("E48").Value = '  +0.83%  '

$ws.Range("D49").Value = "'105.05"
$ws.Range("E49").Value = '  +0.46%  '

$ws.Range("E50").Value = '  +0.01%  '

$ws.Range("D51").Value = "'0.0₆0130"
$ws.Range("E51").Value = '  -5.41%  '

Write-Output "Updated cryptos list"